$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting numeric-looking strings (e.g. "0.574") into numbers.
# Temporarily mark the cell as Text, set the value, then restore the
# default "Normal" style so no visible formatting change is left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "44.185.17"
$ws.Range("E2").Value = "  +3.05%  "

$ws.Range("D3").Value = "2.279.50"
$ws.Range("E3").Value = "  +3.32%  "

$ws.Range("E4").Value = "  -0.30%  "

Set-TextValue $ws.Range("D5") "319.94"
$ws.Range("E5").Value = "  +2.10%  "

Set-TextValue $ws.Range("D6") "103.97"
$ws.Range("E6").Value = "  +6.82%  "

$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("E8").Value = "  -0.23%  "

Set-TextValue $ws.Range("D9") "0.574"
$ws.Range("E9").Value = "  +3.44%  "

Set-TextValue $ws.Range("D10") "39.03"
$ws.Range("E10").Value = "  +7.40%  "

Set-TextValue $ws.Range("D11") "0.0842"
$ws.Range("E11").Value = "  +2.31%  "

$ws.Range("E12").Value = "  +2.81%  "

$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("D14").Value = "2.624.23"
$ws.Range("E14").Value = "  +3.24%  "

Set-TextValue $ws.Range("D15") "0.881"
$ws.Range("E15").Value = "  +2.54%  "

Set-TextValue $ws.Range("D16") "14.61"
$ws.Range("E16").Value = "  +4.19%  "

$ws.Range("D17").Value = "2.279.63"
$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").Value = "44.105.88"
$ws.Range("E18").Value = "  +3.25%  "

Set-TextValue $ws.Range("D19") "14.22"
$ws.Range("E19").Value = "  -2.22%  "

Set-TextValue $ws.Range("D20") "0.0000100"
$ws.Range("E20").Value = "  +4.99%  "

$ws.Range("E21").Value = "  +4.14%  "

Set-TextValue $ws.Range("D22") "66.39"
$ws.Range("E22").Value = "  +2.03%  "

$ws.Range("E23").Value = "  +2.52%  "

Set-TextValue $ws.Range("D24") "238.01"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("E25").Value = "  +4.41%  "

$ws.Range("E26").Value = "  +0.06%  "

Set-TextValue $ws.Range("D27") "10.30"
$ws.Range("E27").Value = "  +3.08%  "

Set-TextValue $ws.Range("D28") "39.15"
$ws.Range("E28").Value = "  +16.62%  "

$ws.Range("E29").Value = "  -0.13%  "

Set-TextValue $ws.Range("D30") "6.58"
$ws.Range("E30").Value = "  +6.77%  "

Set-TextValue $ws.Range("D31") "162.55"
$ws.Range("E31").Value = "  +5.35%  "

Set-TextValue $ws.Range("D32") "20.56"
$ws.Range("E32").Value = "  +0.90%  "

Set-TextValue $ws.Range("D33") "0.0885"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("E34").Value = "  -2.13%  "

Set-TextValue $ws.Range("D35") "2.10"
$ws.Range("E35").Value = "  +6.68%  "

Set-TextValue $ws.Range("D36") "3.28"
$ws.Range("E36").Value = "  +3.69%  "

Set-TextValue $ws.Range("D37") "0.112"
$ws.Range("E37").Value = "  +10.54%  "

$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("E39").Value = "  +3.28%  "

Set-TextValue $ws.Range("D40") "3.90"
$ws.Range("E40").Value = "  +6.28%  "

$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D41") "15.61"
$ws.Range("E41").Value = "  +28.67%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0328"
$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").Value = "1.778.94"
$ws.Range("E44").Value = "  -4.47%  "

$ws.Range("E45").Value = "  +1.42%  "

Set-TextValue $ws.Range("D46") "85.49"
$ws.Range("E46").Value = "  -3.48%  "

Set-TextValue $ws.Range("D47") "5.38"
$ws.Range("E47").Value = "  -0.23%  "

Set-TextValue $ws.Range("D48") "8.93"
$ws.Range("E48").Value = "  +3.64%  "

Set-TextValue $ws.Range("D49") "59.96"
$ws.Range("E49").Value = "  +0.36%  "

Set-TextValue $ws.Range("D50") "75.34"
$ws.Range("E50").Value = "  -0.21%  "

Set-TextValue $ws.Range("D51") "104.76"
$ws.Range("E51").Value = "  +4.36%  "

